$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 23:34"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6821533
$ws.Range("C4").Value = 33386
$ws.Range("D4").Value = 4098391
$ws.Range("E4").Value = 2522049
$ws.Range("G4").Value = 896
$ws.Range("H4").Value = 201093

# Row 6 - Brasil
$ws.Range("B6").Value = 4419083
$ws.Range("C6").Value = 34784
$ws.Range("E6").Value = 613849
$ws.Range("G6").Value = 899
$ws.Range("H6").Value = 134106

# Row 27 - Israel
$ws.Range("B27").Value = 170465
$ws.Range("C27").Value = 6063
$ws.Range("D27").Value = 123219
$ws.Range("E27").Value = 46081
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = 1165

# Row 83 - Costa de Marfil
$ws.Range("B83").Value = 19132
$ws.Range("C83").Value = 32
$ws.Range("D83").Value = 18289
$ws.Range("E83").Value = 723

# Row 167 - Republica del Chad
$ws.Range("B167").Value = 1090
$ws.Range("C167").Value = 3
$ws.Range("D167").Value = 960
$ws.Range("E167").Value = 49
